# Add Figure (from Pandoc 3)
#
# Brings the style sheet up to date with Pandoc 3's reference.docx:
#   - splits the old "Abstract" style into a new "Abstract Title" style
#     (next style = Abstract) plus a retuned "Abstract" style
#   - adds a new "Footnote Block Text" style (pairs with Block Text,
#     but based on Footnote Text)
#   - gives the ImportTok / BuiltInTok syntax-highlighting character
#     styles their colours (they were previously empty placeholders)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. New "Abstract Title" style, shown before the abstract body and
#    flowing into the (retuned) "Abstract" style.
# ---------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("AbstractTitle", 1)
$abstractTitle.NameLocal = "Abstract Title"
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$atPf = $abstractTitle.ParagraphFormat
$atPf.KeepWithNext = $true
$atPf.KeepTogether = $true
$atPf.SpaceBefore = 15
$atPf.SpaceAfter = 0
$atPf.Alignment = 1

$atFont = $abstractTitle.Font
$atFont.Bold = $true
$atFont.Color = 9067060
$atFont.Size = 10
$atFont.SizeBi = 10

# ---------------------------------------------------------------------
# 2. Retune the existing "Abstract" style: space-before drops from
#    300 twips (15pt) to 100 twips (5pt); space-after stays 300 twips.
# ---------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# ---------------------------------------------------------------------
# 3. New "Footnote Block Text" style, mirroring Block Text but based
#    on Footnote Text.
# ---------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("FootnoteBlockText", 1)
$footnoteBlockText.NameLocal = "Footnote Block Text"
$footnoteBlockText.BaseStyle = $d.Styles("FootnoteText")
$footnoteBlockText.NextParagraphStyle = $d.Styles("FootnoteText")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$fbtPf = $footnoteBlockText.ParagraphFormat
$fbtPf.SpaceBefore = 5
$fbtPf.SpaceAfter = 5
$fbtPf.FirstLineIndent = 0
$fbtPf.LeftIndent = 24
$fbtPf.RightIndent = 24

# ---------------------------------------------------------------------
# 4. Colour the two syntax-highlighting character styles that were
#    still empty placeholders.
# ---------------------------------------------------------------------
$importTok = $d.Styles("ImportTok")
$importTok.Font.Bold = $true
$importTok.Font.Color = 32768

$builtInTok = $d.Styles("BuiltInTok")
$builtInTok.Font.Color = 32768

Write-Output "styles updated"
